$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on each D (Price) and E (Volume) cell being updated
# so Excel does not auto-convert numeric-looking strings into numbers.
# (NumberFormat must be set per-cell; applying it to a multi-area range only
# affects the first area.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = '37.482.77'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '2.064.93'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '231.48'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '57.16'
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("D9").Value = '0.387'
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").Value = '0.0776'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '14.81'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '2.372.93'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '20.73'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '0.758'
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").Value = '5.28'
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '2.065.63'
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").Value = '37.392.65'
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").Value = '70.33'
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").Value = '5.90'
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("D22").Value = '227.13'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("D26").Value = '9.54'
$ws.Range("E26").Value = '  +4.66%  '
$ws.Range("D27").Value = '168.76'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '0.131'
$ws.Range("E28").Value = '  -4.03%  '
$ws.Range("D29").Value = '19.27'
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").Value = '1.36'
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("D31").Value = '0.122'
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").Value = '4.56'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").Value = '0.0626'
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  -2.49%  '
$ws.Range("D35").Value = '2.45'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").Value = '3.29'
$ws.Range("E37").Value = '  -3.68%  '
$ws.Range("D39").Value = '5.25'
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").Value = '0.0228'
$ws.Range("E40").Value = '  +6.06%  '
$ws.Range("D41").Value = '98.83'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = '0.0953'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").Value = '2.89'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +3.41%  '
$ws.Range("D45").Value = '1.467.76'
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").Value = '16.41'
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.06'
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '1.03'
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("D49").Value = '7.18'
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("D51").Value = '2.257.21'
$ws.Range("E51").Value = '  -0.51%  '
